# Auto-generated edit script: updates leve-profit calculation values
# across multiple crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 360
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
$ws.Range("H8").Value = 475
$ws.Range("I8").Value = 50.5
$ws.Range("J8").Value = 899.5
$ws.Range("K8").Value = 151.5
$ws.Range("L8").Value = 2698.5
$ws.Range("M8").Value = -12.5
$ws.Range("N8").Value = -2976.5
$ws.Range("H9").Value = 125106.75
$ws.Range("I9").Value = 250063.5
$ws.Range("K9").Value = 250063.5
$ws.Range("M9").Value = -249894.5
$ws.Range("H10").Value = 7500
$ws.Range("J10").Value = 7500
$ws.Range("L10").Value = 7500
$ws.Range("N10").Value = -8086
$ws.Range("H14").Value = 2000
$ws.Range("I14").Value = 2000
$ws.Range("K14").Value = 2000
$ws.Range("M14").Value = -1809
$ws.Range("H16").Value = 1750
$ws.Range("I16").Value = 1750
$ws.Range("K16").Value = 1750
$ws.Range("M16").Value = -1520
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H38").Value = 2230.7144
$ws.Range("I38").Value = 823
$ws.Range("J38").Value = 5750
$ws.Range("K38").Value = 2469
$ws.Range("L38").Value = 17250
$ws.Range("M38").Value = -2097
$ws.Range("N38").Value = -17994
$ws.Range("H39").Value = 271.2
$ws.Range("I39").Value = 89
$ws.Range("K39").Value = 267
$ws.Range("M39").Value = 29
$ws.Range("H48").Value = 1000000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H52").Value = 4750
$ws.Range("I52").Value = 4750
$ws.Range("K52").Value = 14250
$ws.Range("M52").Value = -14090
$ws.Range("H56").Value = 1000000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H58").Value = 995
$ws.Range("I58").Value = 742.5
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 2227.5
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -2077.5
$ws.Range("N58").Value = -4800
$ws.Range("H70").Value = 1868.8
$ws.Range("I70").Value = 1562.6666
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 4687.9998
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -4417.9998
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 1868.8
$ws.Range("I73").Value = 1562.6666
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 4687.9998
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -3751.9998
$ws.Range("N73").Value = -7872
$ws.Range("H132").Value = 2223.9285
$ws.Range("I132").Value = 1624.1818
$ws.Range("K132").Value = 4872.5454
$ws.Range("M132").Value = -2342.5454

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5807.8887
$ws.Range("I32").Value = 5742.846
$ws.Range("K32").Value = 5742.846
$ws.Range("M32").Value = -5455.846
$ws.Range("H97").Value = 931.5
$ws.Range("I97").Value = 914.375
$ws.Range("K97").Value = 914.375
$ws.Range("M97").Value = -418.375
$ws.Range("H102").Value = 2614.3333
$ws.Range("I102").Value = 2172.8333
$ws.Range("J102").Value = 3497.3333
$ws.Range("K102").Value = 2172.8333
$ws.Range("L102").Value = 3497.3333
$ws.Range("M102").Value = -550.8332999999998
$ws.Range("N102").Value = -6741.3333
$ws.Range("H132").Value = 4771.143
$ws.Range("I132").Value = 4975
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 14925
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -12395
$ws.Range("N132").Value = -18558.0005

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 675
$ws.Range("I107").Value = 650
$ws.Range("K107").Value = 650
$ws.Range("M107").Value = 1270

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1116
$ws.Range("I31").Value = 1107.6
$ws.Range("K31").Value = 1107.6
$ws.Range("M31").Value = -812.5999999999999
$ws.Range("H34").Value = 1116
$ws.Range("I34").Value = 1107.6
$ws.Range("K34").Value = 1107.6
$ws.Range("M34").Value = -905.5999999999999
$ws.Range("H94").Value = 2133
$ws.Range("I94").Value = 2699.5
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 2699.5
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -2248.5
$ws.Range("N94").Value = -1902

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2500
$ws.Range("J4").Value = 2500
$ws.Range("L4").Value = 7500
$ws.Range("N4").Value = -7724
$ws.Range("H11").Value = 21231.6
$ws.Range("I11").Value = 160
$ws.Range("J11").Value = 26499.5
$ws.Range("K11").Value = 480
$ws.Range("L11").Value = 79498.5
$ws.Range("M11").Value = -340
$ws.Range("N11").Value = -79778.5
$ws.Range("H26").Value = 350
$ws.Range("I26").Value = 200
$ws.Range("K26").Value = 600
$ws.Range("M26").Value = -312
$ws.Range("H46").Value = 499
$ws.Range("J46").Value = 499
$ws.Range("L46").Value = 1497
$ws.Range("N46").Value = -1679
$ws.Range("H75").Value = 667.5
$ws.Range("I75").Value = 513
$ws.Range("J75").Value = 719
$ws.Range("K75").Value = 1539
$ws.Range("L75").Value = 2157
$ws.Range("N75").Value = -4153
$ws.Range("M75").Value = -541
$ws.Range("H78").Value = 667.5
$ws.Range("I78").Value = 513
$ws.Range("J78").Value = 719
$ws.Range("K78").Value = 4617
$ws.Range("L78").Value = 6471
$ws.Range("N78").Value = -16455
$ws.Range("M78").Value = 375
$ws.Range("H92").Value = 698.8
$ws.Range("I92").Value = 498.33334
$ws.Range("K92").Value = 1495.00002
$ws.Range("M92").Value = -247.0000199999999
$ws.Range("H114").Value = 1233
$ws.Range("I114").Value = 1549.5
$ws.Range("K114").Value = 4648.5
$ws.Range("M114").Value = -1394.5
$ws.Range("H128").Value = 199995.5
$ws.Range("I128").Value = 199995.5
$ws.Range("K128").Value = 599986.5
$ws.Range("M128").Value = -595006.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1842.8572
$ws.Range("I107").Value = 1380.2
$ws.Range("K107").Value = 1380.2
$ws.Range("M107").Value = 539.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4999
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 9893.25
$ws.Range("I22").Value = 15139
$ws.Range("K22").Value = 15139
$ws.Range("M22").Value = -14844
$ws.Range("H27").Value = 9893.25
$ws.Range("I27").Value = 15139
$ws.Range("K27").Value = 15139
$ws.Range("M27").Value = -15032
$ws.Range("H46").Value = 2916.5
$ws.Range("I46").Value = 2916.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2916.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2728.5
$ws.Range("N46").ClearContents()
$ws.Range("H136").Value = 2092.2
$ws.Range("I136").Value = 2154.6667
$ws.Range("K136").Value = 6464.000100000001
$ws.Range("M136").Value = -3914.000100000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 508373.34
$ws.Range("I2").Value = 162863.61
$ws.Range("K2").Value = 162863.61
$ws.Range("M2").Value = -162751.61
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140
